$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (Schneider): Last_update/Next_update values & number-format swap ---
# Last_update (C4) becomes a full datetime, Next_update (D4) becomes a plain date
$ws.Range("C4").Value2 = 44901
$ws.Range("C4").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Range("D4").Value2 = 45266
$ws.Range("D4").NumberFormat = "yyyy\-mm\-dd"

# --- Row 5 now holds the AVM entry (previously on row 6), with Intervall fixed to 0 ---
$ws.Range("A5").Value2 = "AVM"
$ws.Range("B5").Value2 = 0
$ws.Range("C5").Value2 = 44902
$ws.Range("C5").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Range("D5").Value2 = 44902
$ws.Range("D5").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("E5").Value2 = "AVMScraper"

# --- Row 6 now holds the Synology entry (previously on row 5) ---
$ws.Range("A6").Value2 = "Synology"
$ws.Range("B6").Value2 = 0
$ws.Range("C6").Value2 = 44902
$ws.Range("C6").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("D6").Value2 = 44902
$ws.Range("D6").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Range("E6").Value2 = "SynologyScraper"

# --- Cosmetic: widen the Next_update column to fit its contents, move the active selection ---
$ws.Columns.Item(4).ColumnWidth = 17.7265625
$ws.Range("F10").Select() | Out-Null
